$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '60.431.48'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -3.54%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.960.71'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -3.02%  '
$ws.Range('E4').Value = '  -0.03%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '522.50'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -2.02%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '129.25'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -1.98%  '
$ws.Range('E7').Value = '  -0.06%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '2.953.39'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -3.00%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.485'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.147'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -3.76%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '6.08'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -1.11%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.434'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -3.49%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.0000217'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -3.23%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '32.95'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -2.68%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '3.432.44'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('E16').Value = '  -0.28%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '60.341.18'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -3.51%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '2.957.03'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -3.00%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '6.43'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -1.66%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '452.47'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -5.48%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '12.91'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -1.87%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '0.662'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -3.94%  '
$ws.Range('E23').Value = '  -4.06%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '77.50'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -1.33%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '11.63'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -2.79%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -0.34%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.61'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -2.09%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '7.57'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -5.37%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +0.21%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '1.13'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +2.62%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '24.85'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -3.35%  '
$ws.Range('E32').Value = '  -1.45%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '54.56'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -3.51%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '2.22'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -5.57%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '5.27'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -1.06%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '5.71'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -2.76%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '447.80'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -3.84%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '3.156.16'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +2.27%  '
$ws.Range('E39').Value = '  -2.30%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.0374'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -3.96%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.115'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +1.07%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '7.93'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('E43').Value = '  -7.64%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('E45').Value = '  -3.22%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '25.05'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +3.50%  '
$ws.Range('E47').Value = '  -0.38%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '116.89'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -2.96%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.92'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -4.56%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0491'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -10.55%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.25'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +9.54%  '
